$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 77396.62
$ws.Range("I41").Value = 539.0909
$ws.Range("K41").Value = 539.0909
$ws.Range("M41").Value = -99.09090000000003
$ws.Range("H80").Value = 5682390
$ws.Range("I80").Value = 11363974
$ws.Range("J80").Value = 805.63635
$ws.Range("K80").Value = 34091922
$ws.Range("L80").Value = 2416.90905
$ws.Range("M80").Value = -34090924
$ws.Range("N80").Value = -4412.90905
$ws.Range("H83").Value = 5682390
$ws.Range("I83").Value = 11363974
$ws.Range("J83").Value = 805.63635
$ws.Range("K83").Value = 102275766
$ws.Range("L83").Value = 7250.72715
$ws.Range("M83").Value = -102270774
$ws.Range("N83").Value = -17234.72715
$ws.Range("H112").Value = 981.28125
$ws.Range("I112").Value = 1083.3334
$ws.Range("J112").Value = 970.7241
$ws.Range("K112").Value = 3250.0002
$ws.Range("L112").Value = 2912.1723
$ws.Range("M112").Value = -2142.0002
$ws.Range("N112").Value = -5128.1723

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1151.96
$ws.Range("I32").Value = 780.0213
$ws.Range("J32").Value = 6979
$ws.Range("K32").Value = 780.0213
$ws.Range("L32").Value = 6979
$ws.Range("M32").Value = -493.0213
$ws.Range("N32").Value = -7553
$ws.Range("H63").Value = 7500.5
$ws.Range("I63").Value = 6833
$ws.Range("K63").Value = 6833
$ws.Range("M63").Value = -6147
$ws.Range("H66").Value = 7500.5
$ws.Range("I66").Value = 6833
$ws.Range("K66").Value = 34165
$ws.Range("M66").Value = -30733
$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800
$ws.Range("H132").Value = 1865.2458
$ws.Range("I132").Value = 1603.1455
$ws.Range("K132").Value = 4809.4365
$ws.Range("M132").Value = -2279.4365

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 251
$ws.Range("I22").Value = 245
$ws.Range("K22").Value = 245
$ws.Range("M22").Value = -72
$ws.Range("H86").Value = 1569.3077
$ws.Range("J86").Value = 1701.2
$ws.Range("L86").Value = 1701.2
$ws.Range("N86").Value = -3947.2
$ws.Range("H89").Value = 1569.3077
$ws.Range("J89").Value = 1701.2
$ws.Range("L89").Value = 8506
$ws.Range("N89").Value = -19738
$ws.Range("H134").Value = 2291
$ws.Range("I134").Value = 1614.3677
$ws.Range("J134").Value = 3934.25
$ws.Range("K134").Value = 4843.1031
$ws.Range("L134").Value = 11802.75
$ws.Range("M134").Value = -2308.1031
$ws.Range("N134").Value = -16872.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4239.1787
$ws.Range("I31").Value = 1784.4736
$ws.Range("K31").Value = 1784.4736
$ws.Range("M31").Value = -1489.4736
$ws.Range("H34").Value = 4239.1787
$ws.Range("I34").Value = 1784.4736
$ws.Range("K34").Value = 1784.4736
$ws.Range("M34").Value = -1582.4736
$ws.Range("H50").Value = 29285.143
$ws.Range("J50").Value = 31499.5
$ws.Range("L50").Value = 31499.5
$ws.Range("N50").Value = -32749.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 15088
$ws.Range("I3").Value = 756
$ws.Range("K3").Value = 2268
$ws.Range("M3").Value = -2156
$ws.Range("H107").Value = 1416.6666
$ws.Range("I107").Value = 750
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 2250
$ws.Range("L107").Value = 5250
$ws.Range("M107").Value = -330
$ws.Range("N107").Value = -9090
$ws.Range("H121").Value = 25000940
$ws.Range("I121").Value = 50000268
$ws.Range("J121").Value = 1613
$ws.Range("K121").Value = 150000804
$ws.Range("L121").Value = 4839
$ws.Range("M121").Value = -149999494
$ws.Range("N121").Value = -7459
$ws.Range("H139").Value = 169318
$ws.Range("I139").Value = 169318
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 507954
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -502814

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 21942.857
$ws.Range("J57").Value = 21942.857
$ws.Range("L57").Value = 21942.857
$ws.Range("N57").Value = -23582.857
$ws.Range("H70").Value = 6346.6665
$ws.Range("I70").Value = 5116
$ws.Range("K70").Value = 5116
$ws.Range("M70").Value = -4846
$ws.Range("H73").Value = 6346.6665
$ws.Range("I73").Value = 5116
$ws.Range("K73").Value = 5116
$ws.Range("M73").Value = -4180
$ws.Range("H80").Value = 281712.94
$ws.Range("I80").Value = 558528.5600000001
$ws.Range("J80").Value = 4897.3335
$ws.Range("K80").Value = 558528.5600000001
$ws.Range("L80").Value = 4897.3335
$ws.Range("M80").Value = -557530.5600000001
$ws.Range("N80").Value = -6893.3335
$ws.Range("H83").Value = 281712.94
$ws.Range("I83").Value = 558528.5600000001
$ws.Range("J83").Value = 4897.3335
$ws.Range("K83").Value = 2792642.8
$ws.Range("L83").Value = 24486.6675
$ws.Range("M83").Value = -2787650.8
$ws.Range("N83").Value = -34470.6675
$ws.Range("H122").Value = 2539.0557
$ws.Range("I122").Value = 1951.2858
$ws.Range("K122").Value = 5853.857400000001
$ws.Range("M122").Value = -3403.857400000001
$ws.Range("H132").Value = 27787786
$ws.Range("I132").Value = 33340170
$ws.Range("K132").Value = 100020510
$ws.Range("M132").Value = -100017980

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45224
$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45340
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H82").Value = 906.8570999999999
$ws.Range("I82").Value = 787.25
$ws.Range("K82").Value = 787.25
$ws.Range("M82").Value = -426.25
$ws.Range("H85").Value = 906.8570999999999
$ws.Range("I85").Value = 787.25
$ws.Range("K85").Value = 787.25
$ws.Range("M85").Value = 460.75
$ws.Range("H100").Value = 10950
$ws.Range("I100").Value = 2800
$ws.Range("J100").Value = 27250
$ws.Range("K100").Value = 2800
$ws.Range("L100").Value = 27250
$ws.Range("M100").Value = -2259
$ws.Range("N100").Value = -28332
$ws.Range("H122").Value = 5598.5
$ws.Range("I122").Value = 4585.25
$ws.Range("K122").Value = 13755.75
$ws.Range("M122").Value = -11305.75
$ws.Range("H132").Value = 2162
$ws.Range("I132").Value = 2032.6552
$ws.Range("K132").Value = 6097.9656
$ws.Range("M132").Value = -3567.9656

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 99999.5
$ws.Range("I29").Value = 99999.5
$ws.Range("K29").Value = 99999.5
$ws.Range("M29").Value = -99709.5
$ws.Range("H45").Value = 10789
$ws.Range("J45").Value = 11504.6
$ws.Range("L45").Value = 11504.6
$ws.Range("N45").Value = -12486.6
$ws.Range("H107").Value = 2625.3547
$ws.Range("I107").Value = 1575.75
$ws.Range("K107").Value = 4727.25
$ws.Range("M107").Value = -2807.25
